$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 - Dion van Huyssteen
$ws.Range("B24").Value = "won"
$ws.Range("C24").Value = "Digitum Dei"
$ws.Range("F24").Value = "rushed for the middle and then seemed to get confused.  Maybe a poor implementation of a search tree?"

# Row 8 - Isaac Lundall
$ws.Range("B8").Value = "won"
$ws.Range("F8").Value = "one worm rushed for health packs, and then digging commenced"
$ws.Range("C8").Value = "Ikabot"

# Row 5 - Justin Wernick
$ws.Range("B5").Value = "lost"
$ws.Range("C5").Value = "Steam powered worm"
$ws.Range("F5").Value = "lost by KO.  He was doing some fancy dodging in the skirmishes! But I still won the first one tho"

# Row 17 - Marvin Thobejane
$ws.Range("B17").Value = "won"
$ws.Range("C17").Value = "Marvijo"
$ws.Range("F17").Value = "rushed for the middle and then just hung around there?  I went on digging merrily and won by points"

# Update the view's scroll position and selection to reflect where the
# author ended up working (scrolled down so row 7 is the top visible row,
# with B18 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("B18").Select()
